# design_input_background.xlsx update:
# - Add correlated PARAM5 ~ TruncatedNormal(3,1,1,5) and PARAM6 ~ Uniform(0,1)
#   via a new "corr0" correlation sheet, inserted before "corr1".
# - Point the designinput rows for PARAM5 / PARAM6 (rows 9 & 10) at the new
#   "corr0" correlation sheet, and bump numreal for PARAM5 (B9) to 500.

$wb = $excel.ActiveWorkbook

# Insert the new correlation sheet right before "corr1" so the tab order
# becomes: general_input, designinput, defaultvalues, corr0, corr1,
# background, background_corr.
$corr1 = $wb.Worksheets.Item("corr1")
$corr0 = $wb.Worksheets.Add($corr1)
$corr0.Name = "corr0"

# 2x2 correlation matrix between PARAM5 and PARAM6 (PARAM6 correlates 0.8
# with PARAM5), following the same layout used by the existing "corr1"
# sheet.
$corr0.Range("B1").Value = "PARAM5"
$corr0.Range("C1").Value = "PARAM6"

$corr0.Range("A2").Value = "PARAM5"
$corr0.Range("B2").Value = 1

$corr0.Range("A3").Value = "PARAM6"
$corr0.Range("B3").Value = 0.8
$corr0.Range("C3").Value = 1

# Update the designinput sheet: PARAM5 (row 9) and PARAM6 (row 10) now use
# the corr0 correlation sheet, and PARAM5's numreal changes from 10 to 500.
$design = $wb.Worksheets.Item("designinput")
$design.Range("B9").Value = 500
$design.Range("O9").Value = "corr0"
$design.Range("O10").Value = "corr0"
$design.Range("B10").Select() | Out-Null

# Leave the new corr0 sheet as the active tab/selection.
$corr0.Activate()
$corr0.Range("C8").Select() | Out-Null
